# Insert a new data row at row 357 (pushing existing rows 357-468 down to 358-469)
# and populate it with the new weekly price record for Acelga.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(357).Insert()

$ws.Cells.Item(357, 1).Value = 10
$ws.Cells.Item(357, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(357, 3).Value = "La Araucanía"
$ws.Cells.Item(357, 4).Value = 44988
$ws.Cells.Item(357, 5).Value = 9
$ws.Cells.Item(357, 6).Value = 100112009
$ws.Cells.Item(357, 7).Value = "Acelga"
$ws.Cells.Item(357, 8).Value = "Sin especificar"
$ws.Cells.Item(357, 9).Value = "Primera"
$ws.Cells.Item(357, 10).Value = 30
$ws.Cells.Item(357, 11).Value = 8000
$ws.Cells.Item(357, 12).Value = 8000
$ws.Cells.Item(357, 13).Value = 8000
$ws.Cells.Item(357, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(357, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(357, 16).Value = 667
$ws.Cells.Item(357, 17).Value = 12
$ws.Cells.Item(357, 18).Value = "Hortaliza"
